# Normalize vignette_type (column A) and evaluator_gender (column G) values
# to title case across all data rows (rows 2-67).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 67 }

$map = @{
    "higher-class" = "Higher-Class"
    "lower-class"  = "Lower-Class"
    "male"         = "Male"
    "female"       = "Female"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($map.ContainsKey($aVal)) {
        $aCell.Value2 = $map[$aVal]
    }

    $gCell = $ws.Cells.Item($r, 7)
    $gVal = $gCell.Value2
    if ($map.ContainsKey($gVal)) {
        $gCell.Value2 = $map[$gVal]
    }
}
